$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap match data (columns F:V) between row 77 and row 78 ---
# Row 77 (Varzim 2 - 1 Sanjoanense)  <->  Row 78 (Lusitania FC 4 - 0 Anadia)
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row77vals = @{}
$row78vals = @{}
foreach ($col in $cols) {
    $row77vals[$col] = $ws.Range($col + "77").Value2
    $row78vals[$col] = $ws.Range($col + "78").Value2
}
foreach ($col in $cols) {
    $ws.Range($col + "77").Value = $row78vals[$col]
    $ws.Range($col + "78").Value = $row77vals[$col]
}

# --- 2. Append new row 91 (new match: Alverca 1 - 0 Oliveira Hospital) ---
# Copy formatting from the last existing data row (90) onto the new row.
$ws.Range("A90:V90").Copy()
$ws.Range("A91:V91").PasteSpecial(-4122)

$ws.Range("A91").Value = 90
$ws.Range("B91").Value = "portugal"
$ws.Range("C91").Value = "liga-3"
$ws.Range("D91").Value = "2023-2024"
$ws.Range("E91").Value = 45233.85416666666
$ws.Range("F91").Value = "Alverca"
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = "Oliveira Hospital"
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1.71
$ws.Range("K91").Value = "02/11/2023 13:42"
$ws.Range("L91").Value = 1.65
$ws.Range("M91").Value = "03/11/2023 20:29"
$ws.Range("N91").Value = 3.64
$ws.Range("O91").Value = "02/11/2023 13:42"
$ws.Range("P91").Value = 3.79
$ws.Range("Q91").Value = "03/11/2023 20:29"
$ws.Range("R91").Value = 4.98
$ws.Range("S91").Value = "02/11/2023 13:42"
$ws.Range("T91").Value = 5.65
$ws.Range("U91").Value = "03/11/2023 20:29"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/portugal/liga-3/alverca-oliveira-hospital/vVkOh1YP/"
